$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing "placeholder" text cells (shared strings change in place) ---
# Row 7: "as" -> "ssdf"
$ws.Range("C7").Value = "ssdf"
$ws.Range("J7").Value = "ssdf"

# Row 8: "asd" -> "sadfsdg"
$ws.Range("C8").Value = "sadfsdg"
$ws.Range("J8").Value = "sadfsdg"

# Row 9: "asd" -> "dafsdf" (new distinct text)
$ws.Range("C9").Value = "dafsdf"
$ws.Range("J9").Value = "dafsdf"

# --- Add the new "end of sheet" block: signatures / approval section ---

# Row 26: "Zespół Orzekający:" (merged A:B and H:I) plus "Zatwierdzam" at F26/M26
$ws.Range("A26").Value = "Zespół Orzekający:"
$ws.Range("H26").Value = "Zespół Orzekający:"
$ws.Range("F26").Value = "Zatwierdzam"
$ws.Range("M26").Value = "Zatwierdzam"

# Row 29: "1 ................" (merged A:B and H:I)
$ws.Range("A29").Value = "1 ................"
$ws.Range("H29").Value = "1 ................"

# Row 33: "2 ................" (merged A:B and H:I)
$ws.Range("A33").Value = "2 ................"
$ws.Range("H33").Value = "2 ................"

# Left-align the merged label cells (matches the new cellXfs entry: borderId 0 + horizontal left)
$ws.Range("A26:B26").HorizontalAlignment = -4131
$ws.Range("H26:I26").HorizontalAlignment = -4131
$ws.Range("A29:B29").HorizontalAlignment = -4131
$ws.Range("H29:I29").HorizontalAlignment = -4131
$ws.Range("A33:B33").HorizontalAlignment = -4131
$ws.Range("H33:I33").HorizontalAlignment = -4131

# Merge the label cells (order matches the canonical mergeCells list)
$ws.Range("A26:B26").Merge()
$ws.Range("A29:B29").Merge()
$ws.Range("A33:B33").Merge()
$ws.Range("H26:I26").Merge()
$ws.Range("H29:I29").Merge()
$ws.Range("H33:I33").Merge()
